# Apply "write_partial_takes and write_participation" edit
# - Fill row 7 (E7:BB7) on Pag2..Pag9 with sequential "take" numbers
#   continuing from the previous sheet's range (1-50, 51-100, ... 401-450)
# - Update the last active-cell selection remembered on each sheet

$wb = $excel.ActiveWorkbook

function Fill-Row7 {
    param(
        [string]$SheetName,
        [int]$StartValue
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
              "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU","AV","AW","AX","AY","AZ", `
              "BA","BB")
    $value = $StartValue
    foreach ($col in $cols) {
        $ws.Range($col + "7").Value = $value
        $value = $value + 1
    }
}

function Set-Selection {
    param(
        [string]$SheetName,
        [string]$CellAddress
    )
    $ws = $wb.Worksheets.Item($SheetName)
    $ws.Activate()
    $ws.Range($CellAddress).Select()
}

# Pag1: only the remembered selection changes (no row7 data on this sheet)
Set-Selection "Pag1" "R17"

# Pag2: take numbers 51-100
Fill-Row7 "Pag2" 51
Set-Selection "Pag2" "BB8"

# Pag3: take numbers 101-150
Fill-Row7 "Pag3" 101
Set-Selection "Pag3" "BB7"

# Pag4: take numbers 151-200
Fill-Row7 "Pag4" 151
Set-Selection "Pag4" "BB8"

# Pag5: take numbers 201-250
Fill-Row7 "Pag5" 201
Set-Selection "Pag5" "BB7"

# Pag6: take numbers 251-300
Fill-Row7 "Pag6" 251
Set-Selection "Pag6" "BB7"

# Pag7: take numbers 301-350
Fill-Row7 "Pag7" 301
Set-Selection "Pag7" "AX18"

# Pag8: take numbers 351-400
Fill-Row7 "Pag8" 351
Set-Selection "Pag8" "AD19"

# Pag9: take numbers 401-450
Fill-Row7 "Pag9" 401
Set-Selection "Pag9" "S8"

# Leave the workbook focused back on the first sheet (it is the tabSelected sheet)
$wb.Worksheets.Item("Pag1").Activate()
